$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '[Daniel J%Escobar%Daniel.escobar@pennmedicine.upenn.edu%1,                  Maria%Lanzi%NULL%1,                  Pouné%Saberi%NULL%1,                  Ruby%Love%NULL%1,                  Darren R%Linkin%NULL%1,                  John J%Kelly%NULL%1,                  Darshana%Jhala%NULL%1,                  Valerianna%Amorosa%NULL%1,                  Mary%Hofmann%NULL%1,                  Jeffrey B%Doyon%NULL%1]'
$ws.Range("I2").Value = '_PMC_CROSSREF'
$ws.Range("E3").Value = '[James L.%Rudolph%NULL%1,                  Christopher W.%Halladay%NULL%1,                  Malisa%Barber%NULL%1,                  Kevin W.%McConeghy%NULL%1,                  Vince%Mor%NULL%1,                  Aman%Nanda%NULL%1,                  Stefan%Gravenstein%NULL%1]'
$ws.Range("I3").Value = '_PMC_elsevier_CROSSREF'
$ws.Range("E4").Value = '[Paula%Eckardt%NULL%1,                  Rachel%Guran%NULL%1,                  Jon%Hennemyre%NULL%1,                  Roshan%Arikupurathu%NULL%1,                  Julie%Poveda%NULL%1,                  Nancimae%Miller%NULL%1,                  Randy%Katz%NULL%1,                  Judith%Frum%NULL%1]'
$ws.Range("I4").Value = '_PMC_elsevier_CROSSREF'
$ws.Range("E6").Value = '[Carson%Telford%NULL%1,                  Udodirim%Onwubiko%NULL%1,                  David%Holland%NULL%1,                  Kim%Turner%NULL%1,                  Juliana%Prieto%NULL%1,                  Sasha%Smith%NULL%1,                  Jane%Yoon%NULL%1,                  Wecheeta%Brown%NULL%1,                  Allison%Chamberlain%NULL%1,                  Neel%Gandhi%NULL%1,                  Shamimul%Khan%NULL%1,                  Steve%Williams%NULL%1,                  Fazle%Khan%NULL%1,                  Sarita%Shah%NULL%1,                   C. T.%Telford%null%1,                   U.% Onwubiko%null%1,                   D.% Holland%null%1,                   K.% Turner%null%1,                   J.% Prieto%null%1,                   S.% Smith%null%1,                   J.% Yoon%null%1,                   W.% Brown%null%1,                   A.% Chamberlain%null%1,                   N.% Gandhi%null%1,                   S.% Khan%null%1,                   S.% Williams%null%1,                   F.% Khan%null%1,                   S. % Shah%null%1,                  C. T.%Telford%null%1,                  U.% Onwubiko%null%1,                  D.% Holland%null%1,                  K.% Turner%null%1,                  J.% Prieto%null%1,                  S.% Smith%null%1,                  J.% Yoon%null%1,                  W.% Brown%null%1,                  A.% Chamberlain%null%1,                  N.% Gandhi%null%1,                  S.% Khan%null%1,                  S.% Williams%null%1,                  F.% Khan%null%1,                  S. % Shah%null%1]'
$ws.Range("I6").Value = '_PMC_CROSSREF'
$ws.Range("E7").Value = '[Sandra M.%Shi%NULL%1,                  Innokentiy%Bakaev%NULL%1,                  Helen%Chen%NULL%1,                  Thomas G.%Travison%NULL%1,                  Sarah D.%Berry%NULL%1]'
$ws.Range("I7").Value = '_PMC_elsevier_CROSSREF'
$ws.Range("E8").Value = '[Temet M.%McMichael%NULL%0,                  Dustin W.%Currie%NULL%0,                  Shauna%Clark%NULL%0,                  Sargis%Pogosjans%NULL%0,                  Meagan%Kay%NULL%0,                  Noah G.%Schwartz%NULL%0,                  James%Lewis%NULL%0,                  Atar%Baer%NULL%0,                  Vance%Kawakami%NULL%0,                  Margaret D.%Lukoff%NULL%0,                  Jessica%Ferro%NULL%0,                  Claire%Brostrom-Smith%NULL%0,                  Thomas D.%Rea%NULL%1,                  Michael R.%Sayre%NULL%1,                  Francis X.%Riedo%NULL%0,                  Denny%Russell%NULL%0,                  Brian%Hiatt%NULL%0,                  Patricia%Montgomery%NULL%0,                  Agam K.%Rao%NULL%0,                  Eric J.%Chow%NULL%0,                  Farrell%Tobolowsky%NULL%0,                  Michael J.%Hughes%NULL%1,                  Ana C.%Bardossy%NULL%0,                  Lisa P.%Oakley%NULL%0,                  Jesica R.%Jacobs%NULL%0,                  Nimalie D.%Stone%NULL%1,                  Sujan C.%Reddy%NULL%0,                  John A.%Jernigan%NULL%0,                  Margaret A.%Honein%NULL%0,                  Thomas A.%Clark%NULL%0,                  Jeffrey S.%Duchin%NULL%0]'
$ws.Range("I8").Value = '_PMC_CROSSREF'
$ws.Range("I9").Value = '_CROSSREF'
$ws.Range("E10").Value = '[Nathan M.%Stall%nathan.stall@sinaihealth.ca%1,                  Carolyn%Farquharson%NULL%2,                  Carolyn%Farquharson%NULL%0,                  Chris%Fan‐Lun%NULL%1,                  Lesley%Wiesenfeld%NULL%1,                  Carla A.%Loftus%NULL%1,                  Dylan%Kain%NULL%1,                  Jennie%Johnstone%NULL%1,                  Liz%McCreight%NULL%1,                  Russell D.%Goldman%NULL%1,                  Ramona%Mahtani%NULL%1]'
$ws.Range("I10").Value = '_PMC_CROSSREF'
$ws.Range("E11").Value = '[ Nathan M.%Stall%null%1,                   Aaron%Jones%null%1,                   Kevin A.%Brown%null%1,                   Paula A.%Rochon%null%1,                   Andrew P.%Costa%null%1,                  Nathan M.%Stall%null%1,                  Aaron%Jones%null%1,                  Kevin A.%Brown%null%1,                  Paula A.%Rochon%null%1,                  Andrew P.%Costa%null%1]'
$ws.Range("I11").Value = '_CROSSREF'
$ws.Range("E12").Value = '[Hubert%Blain%NULL%1,                  Yves%Rolland%NULL%1,                  Edouard%Tuaillon%NULL%1,                  Nadia%Giacosa%NULL%1,                  Mylène%Albrand%NULL%1,                  Audrey%Jaussent%NULL%1,                  Athanase%Benetos%NULL%1,                  Stéphanie%Miot%NULL%1,                  Jean%Bousquet%NULL%1]'
$ws.Range("I12").Value = '_PMC_elsevier_CROSSREF'
$ws.Range("C13").Value = '"Coronavirus Disease 2019 Outcomes in French Nursing Homes That Implemented Staff Confinement With Residents"'
$ws.Range("E13").Value = '[Jo\u00ebl%Belmin%xref no email%1, Nathavy%Um-Din%xref no email%1, Cristiano%Donadio%xref no email%1, Maurizio%Magri%xref no email%1, Quoc Duy%Nghiem%xref no email%1, Bruno%Oquendo%xref no email%1, Sylvie%Pariel%xref no email%1, Carmelo%Lafuente-Lafuente%xref no email%1]'
$ws.Range("F13").Value = '10.1001/jamanetworkopen.2020.17533'
$ws.Range("G13").Value = 'CROSSREF'
$ws.Range("H13").Value = '2023-05-11'
$ws.Range("E14").Value = '[Guillaume%Sacco%NULL%1,                  Gonzague%Foucault%NULL%1,                  Olivier%Briere%NULL%1,                  Cédric%Annweiler%NULL%1]'
$ws.Range("I14").Value = '_PMC_elsevier_CROSSREF'
$ws.Range("E15").Value = '[R.%Guery%NULL%1,                  C.%Delaye%NULL%1,                  N.%Brule%NULL%1,                  V.%Nael%NULL%1,                  L.%Castain%NULL%1,                  F.%Raffi%NULL%1,                  L.%De Decker%NULL%1]'
$ws.Range("I15").Value = '_PMC_elsevier_CROSSREF'
$ws.Range("E16").Value = '[A.%Klein%anke.klein@uke.de%1,                  C.%Edler%NULL%1,                  A.%Fitzek%NULL%1,                  D.%Fröb%NULL%1,                  A.%Heinemann%NULL%1,                  K.%Meißner%NULL%1,                  H.%Mushumba%NULL%1,                  K.%Püschel%NULL%1,                  A. S.%Schröder%NULL%1,                  J. P.%Sperhake%NULL%1,                  F.%Ishorst-Witte%NULL%1,                  M.%Aepfelbacher%NULL%1,                  F.%Heinrich%NULL%1]'
$ws.Range("I16").Value = '_PMC_Springer_CROSSREF'
$ws.Range("E17").Value = '[Buntinx%Frank%coreGivesNoEmail%1,                Claes%Peter%coreGivesNoEmail%1,                Gulikers%Marjo%coreGivesNoEmail%1,                Jan%De Lepeleire%coreGivesNoEmail%1,                Van%der Elst Michael%coreGivesNoEmail%1,                Van%Elslande Jan%coreGivesNoEmail%1,                Van%Ranst Marc%coreGivesNoEmail%1,                Verbakel%Jan%coreGivesNoEmail%1,                Vernneersch%Pieter%coreGivesNoEmail%1]'
$ws.Range("I17").Value = '_CROSSREF'
$ws.Range("I18").Value = '_PMC_CROSSREF'
$ws.Range("E19").Value = '[Liangcong%Hu%xref no email%1,            Wanyuan%Zhen%xref no email%1,            Xudong%Xie%xref no email%2,            Ze%Lin%xref no email%2,            Tiantian%Wang%xref no email%1,            Jing%Liu%xref no email%5,            Hang%Xue%xref no email%2,            Adriana C.%Panayi%xref no email%2,            Ke%Xu%xref no email%1,            Zexi%Ling%xref no email%1,            Xugui%Li%xref no email%1,            Bobin%Mi%xref no email%1,            Wu%Zhou%xref no email%2,            Guohui%Liu%xref no email%2]'
$ws.Range("I19").Value = '_CROSSREF'
$ws.Range("E20").Value = '[Sean P%Kennelly%Sean.Kennelly@tuh.ie%1,                  Adam H%Dyer%dyera@tcd.ie%1,                  Claire%Noonan%claire.noonan@tuh.ie%2,                  Claire%Noonan%claire.noonan@tuh.ie%0,                  Ruth%Martin%ruth.martin1@hse.ie%1,                  Siobhan M%Kennelly%siobhan.kennelly1@hse.ie%1,                  Alan%Martin%alanmartin@beaumont.ie%1,                  Desmond%O’Neill%NULL%1,                  Aoife%Fallon%aoife.fallon1@hotmail.com%1]'
$ws.Range("I20").Value = '_PMC_CROSSREF'
$ws.Range("E21").Value = '[Antonio%Nouvenne%NULL%1,                  Andrea%Ticinesi%NULL%1,                  Alberto%Parise%NULL%1,                  Beatrice%Prati%NULL%1,                  Marcello%Esposito%NULL%1,                  Valentina%Cocchi%NULL%1,                  Emanuele%Crisafulli%NULL%1,                  Annalisa%Volpi%NULL%1,                  Sandra%Rossi%NULL%1,                  Elena Giovanna%Bignami%NULL%1,                  Marco%Baciarello%NULL%1,                  Ettore%Brianti%NULL%1,                  Massimo%Fabi%NULL%1,                  Tiziana%Meschi%NULL%1]'
$ws.Range("I21").Value = '_PMC_elsevier_CROSSREF'
$ws.Range("E22").Value = '[Nicola%Veronese%NULL%1,                  Luca Gino%Sbrogiò%NULL%1,                  Roberto%Valle%NULL%1,                  Laura%Marin%NULL%1,                  Elena%Boscolo Fiore%NULL%1,                  Andrea%Tiozzo%NULL%1]'
$ws.Range("I22").Value = '_PMC_elsevier_CROSSREF'
$ws.Range("I23").Value = '_CROSSREF'
$ws.Range("E24").Value = '[ J. H.%van den Besselaar%null%1,                   R. S.% Sikkema%null%1,                   F. M. H. P. A.% Koene%null%1,                   L. W.% van Buul%null%1,                   B. B.% Oude Munnink%null%1,                   I.% Frenay%null%1,                   R.% te Witt%null%1,                   M. P. G.% Koopmans%null%1,                   C. M. P. M.% Hertogh%null%1,                   B. M. % Buurman%null%1]'
$ws.Range("I24").Value = '_MedBiorxiv_CROSSREF'
$ws.Range("E25").Value = '[Laura W.%van Buul%NULL%1,                  Judith H.%van den Besselaar%NULL%2,                  Judith H.%van den Besselaar%NULL%0,                  Fleur M. H. P. H.%Koene%NULL%1,                  Bianca M.%Buurman%NULL%1,                  Cees M. P. M.%Hertogh%NULL%1,                  Martin%Smalbrugge%NULL%1,                  Jeanine J. S.%Rutten%NULL%1,                  Elke M.%den Boogert%NULL%1,                  Michel D.%Wissing%NULL%1,                  Ariene%Rietveld%NULL%1,                  Mariska W. W.%van Elsakker%NULL%1,                  Marga M. G.%Nonneman%NULL%1,                  Florien%van Eeden%NULL%1,                  Saskia%van de Merwe%NULL%1,                  Sophie L.%Niemansburg%NULL%1,                  Ewout%Fanoy%NULL%1,                  Hinke S.%Bootsma%NULL%1,                  Nicoline%van der Hagen%NULL%1,                  Mariska%Petrignani%NULL%1,                  Jessica Edwards%van Muijen%NULL%1,                  Karolien E. M.%Biesheuvel%NULL%1]'
$ws.Range("I25").Value = '_PMC_CROSSREF'
$ws.Range("E26").Value = '[B\u00e5rd%Reiakvam Kittang%coreGivesNoEmail%1,                Karina%Koller L\u00f8land%coreGivesNoEmail%1,                Kjell%Kr\u00fcger%coreGivesNoEmail%1,                Kristian%Jansen%coreGivesNoEmail%1,                Sabine%Piepenstock Solheim%coreGivesNoEmail%1,                Sebastian%von Hofacker%coreGivesNoEmail%1]'
$ws.Range("I26").Value = '_CROSSREF'
$ws.Range("E27").Value = '[Shin Young%Park%NULL%1,                  Gawon%Choi%NULL%2,                  Gawon%Choi%NULL%0,                  Hyeyoung%Lee%NULL%2,                  Hyeyoung%Lee%NULL%0,                  Na-young%Kim%NULL%2,                  Na-young%Kim%NULL%0,                  Seon-young%Lee%NULL%2,                  Seon-young%Lee%NULL%0,                  Kyungnam%Kim%NULL%2,                  Kyungnam%Kim%NULL%0,                  Soyoung%Shin%NULL%2,                  Soyoung%Shin%NULL%0,                  Eunsu%Jang%NULL%2,                  Eunsu%Jang%NULL%0,                  YoungSin%Moon%NULL%2,                  YoungSin%Moon%NULL%0,                  KwangHwan%Oh%NULL%2,                  KwangHwan%Oh%NULL%0,                  JaeRin%Choi%NULL%2,                  JaeRin%Choi%NULL%0,                  Sangeun%Lee%NULL%2,                  Sangeun%Lee%NULL%0,                  Young-Man%Kim%NULL%2,                  Young-Man%Kim%NULL%0,                  Jieun%Kim%NULL%2,                  Jieun%Kim%NULL%0,                  Seonju%Yi%NULL%2,                  Seonju%Yi%NULL%0,                  Jin%Gwack%NULL%2,                  Jin%Gwack%NULL%0,                  Ok%Park%NULL%2,                  Ok%Park%NULL%0,                  Young Joon%Park%NULL%2,                  Young Joon%Park%NULL%0]'
$ws.Range("I27").Value = '_PMC_CROSSREF'
$ws.Range("E28").Value = '[Rok%Song%NULL%1,                  Hee-Sook%Kim%NULL%2,                  Hee-Sook%Kim%NULL%0,                  Seok-Ju%Yoo%NULL%2,                  Seok-Ju%Yoo%NULL%0,                  Kwan%Lee%NULL%2,                  Kwan%Lee%NULL%0,                  Ji-Hyuk%Park%NULL%2,                  Ji-Hyuk%Park%NULL%0,                  Joon Ho%Jang%NULL%2,                  Joon Ho%Jang%NULL%0,                  Gyoung-Sook%Ahn%NULL%2,                  Gyoung-Sook%Ahn%NULL%0,                  Jun-Nyun%Kim%NULL%2,                  Jun-Nyun%Kim%NULL%0]'
$ws.Range("I28").Value = '_PMC_CROSSREF'
$ws.Range("E29").Value = '[Blanca%Borras-Bermejo%NULL%1,                  Xavier%Martínez-Gómez%NULL%1,                  María Gutierrez%San Miguel%NULL%1,                  Juliana%Esperalba%NULL%1,                  Andrés%Antón%NULL%1,                  Elisabet%Martin%NULL%1,                  Marta%Selvi%NULL%1,                  María José%Abadías%NULL%1,                  Antonio%Román%NULL%1,                  Tomàs%Pumarola%NULL%1,                  Magda%Campins%NULL%1,                  Benito%Almirante%NULL%1]'
$ws.Range("I29").Value = '_PMC_CROSSREF'
$ws.Range("E30").Value = '[M.%Bernabeu-Wittel%NULL%1,                  J.E.%Ternero-Vega%NULL%1,                  P.%Díaz-Jiménez%NULL%1,                  C.%Conde-Guzmán%NULL%1,                  M.D.%Nieto-Martín%NULL%1,                  L.%Moreno-Gaviño%NULL%1,                  J.%Delgado-Cuesta%NULL%1,                  M.%Rincón-Gómez%NULL%1,                  L.%Giménez-Miranda%NULL%1,                  M.D.%Navarro-Amuedo%NULL%1,                  M.M.%Muñoz-García%NULL%1,                  S.%Calzón-Fernández%NULL%1,                  M.%Ollero-Baturone%NULL%1]'
$ws.Range("I30").Value = '_PMC_elsevier_CROSSREF'
$ws.Range("E31").Value = '[Shamez N%Ladhani%shamez.ladhani@phe.gov.uk%1,                  J.Yimmy%Chow%NULL%1,                  Roshni%Janarthanan%NULL%1,                  Jonathan%Fok%NULL%1,                  Emma%Crawley-Boevey%NULL%1,                  Amoolya%Vusirikala%NULL%1,                  Elena%Fernandez%NULL%1,                  Marina Sanchez%Perez%NULL%1,                  Suzanne%Tang%NULL%1,                  Kate%Dun-Campbell%NULL%1,                  Edward Wynne-%Evans%NULL%1,                  Anita%Bell%NULL%1,                  Bharat%Patel%NULL%1,                  Zahin%Amin-Chowdhury%NULL%1,                  Felicity%Aiano%NULL%1,                  Karthik%Paranthaman%NULL%1,                  Thomas%Ma%NULL%1,                  Maria%Saavedra-Campos%NULL%1,                  Richard%Myers%NULL%1,                  Joanna%Ellis%NULL%2,                  Angie%Lackenby%NULL%1,                  Robin%Gopal%NULL%1,                  Monika%Patel%NULL%1,                  Colin%Brown%NULL%1,                  Meera%Chand%NULL%1,                  Kevin%Brown%NULL%1,                  Mary E%Ramsay%NULL%1,                  Susan%Hopkins%NULL%2,                  Nandini%Shetty%NULL%1,                  Maria%Zambon%NULL%2]'
$ws.Range("I31").Value = '_PMC_elsevier_CROSSREF'
$ws.Range("E32").Value = '[N.S.N.%Graham%NULL%1,                  C.%Junghans%NULL%1,                  R.%Downes%NULL%1,                  C.%Sendall%NULL%1,                  H.%Lai%NULL%1,                  A.%McKirdy%NULL%1,                  P.%Elliott%NULL%1,                  R.%Howard%NULL%1,                  D.%Wingfield%NULL%1,                  M.%Priestman%NULL%1,                  M.%Ciechonska%NULL%1,                  L.%Cameron%NULL%1,                  M.%Storch%NULL%1,                  M.A.%Crone%NULL%1,                  P.S.%Freemont%NULL%1,                  P.%Randell%NULL%1,                  R.%McLaren%NULL%1,                  N.%Lang%NULL%1,                  S.%Ladhani%NULL%1,                  F.%Sanderson%NULL%1,                  D.J.%Sharp%NULL%1]'
$ws.Range("I32").Value = '_PMC_elsevier_CROSSREF'
$ws.Range("E33").Value = '[Agnes%Marossy%agnesmarossy@nhs.net%1,                  Stefan%Rakowicz%NULL%1,                  Angela%Bhan%NULL%1,                  Sarah%Noon%NULL%1,                  Amanda%Rees%NULL%1,                  Manjinder%Virk%NULL%1,                  Ayazali%Nazafi%NULL%1,                  Evie%Hay%NULL%1,                  Louise%de Thomasson%NULL%1,                  Christina%Windle%NULL%1,                  Mark%Zuckerman%Mark.Zuckerman@nhs.net%1]'
$ws.Range("I33").Value = '_PMC_CROSSREF'
$ws.Range("E34").Value = '[Laura%Shallcross%NULL%1,                  Danielle%Burke%NULL%1,                  Owen%Abbott%NULL%1,                  Alasdair%Donaldson%NULL%1,                  Gemma%Hallatt%NULL%1,                  Andrew%Hayward%NULL%1,                  Susan%Hopkins%NULL%0,                  Maria%Krutikov%NULL%1,                  Katie%Sharp%NULL%1,                  Leone%Wardman%NULL%1,                  Sapphira%Thorne%NULL%1]'
$ws.Range("I34").Value = '_PMC_elsevier_CROSSREF'
$ws.Range("E35").Value = '[ E.%Smith%null%1,                   C. F.% Aldus%null%1,                   J.% Brainard%null%1,                   S.% Dunham%null%1,                   P. R.% Hunter%null%1,                   N.% Steel%null%1,                   P. % Everden%null%1]'
$ws.Range("I35").Value = '_MedBiorxiv_CROSSREF'
$ws.Range("E36").Value = '[Mahesh C%Patel%mp3@uic.edu%1,                  Lelia H%Chaisson%NULL%1,                  Scott%Borgetti%NULL%1,                  Deborah%Burdsall%NULL%1,                  Rashmi K%Chugh%NULL%1,                  Christopher R%Hoff%NULL%1,                  Elizabeth B%Murphy%NULL%1,                  Emily A%Murskyj%NULL%1,                  Shannon%Wilson%NULL%1,                  Joe%Ramos%NULL%1,                  Lynn%Akker%NULL%1,                  Debra%Bryars%NULL%1,                  Evonda%Thomas-Smith%NULL%1,                  Susan C%Bleasdale%NULL%1,                  Ngozi O%Ezike%NULL%1]'
$ws.Range("I36").Value = '_PMC_CROSSREF'
$ws.Range("E37").Value = '[Helena%Temkin‐Greener%helena_temkin-greener@urmc.rochester.edu%1,                  Wenhan%Guo%NULL%2,                  Wenhan%Guo%NULL%0,                  Yunjiao%Mao%NULL%1,                  Xueya%Cai%NULL%1,                  Yue%Li%NULL%1]'
$ws.Range("I37").Value = '_PMC_CROSSREF'
$ws.Range("E38").Value = '[Ana A.%Weil%anaweil@uw.edu%1,                  Kira L.%Newman%NULL%2,                  Kira L.%Newman%NULL%0,                  Thuan D.%Ong%NULL%1,                  Giana H.%Davidson%NULL%1,                  Jennifer%Logue%NULL%1,                  Elisabeth%Brandstetter%NULL%1,                  Ariana%Magedson%NULL%1,                  Dylan%McDonald%NULL%1,                  Denise J.%McCulloch%NULL%1,                  Santiago%Neme%NULL%1,                  James%Lewis%NULL%0,                  Jeff S.%Duchin%NULL%1,                  Weizhi%Zhong%NULL%1,                  Lea M.%Starita%NULL%1,                  Trevor%Bedford%NULL%1,                  Alison C.%Roxby%NULL%1,                  Helen Y.%Chu%helenchu@uw.edu%2]'
$ws.Range("I38").Value = '_PMC_Springer_CROSSREF'
$ws.Range("E39").Value = '[Guillermo V.%Sanchez%NULL%1,                  Caitlin%Biedron%NULL%1,                  Lauren R.%Fink%NULL%1,                  Kelly M.%Hatfield%NULL%1,                  Jordan Micah F.%Polistico%NULL%1,                  Monica P.%Meyer%NULL%1,                  Rebecca S.%Noe%NULL%1,                  Casey E.%Copen%NULL%1,                  Amanda K.%Lyons%NULL%1,                  Gonzalo%Gonzalez%NULL%1,                  Keith%Kiama%NULL%1,                  Mark%Lebednick%NULL%1,                  Bonnie K.%Czander%NULL%1,                  Amen%Agbonze%NULL%1,                  Aimee R.%Surma%NULL%1,                  Avnish%Sandhu%NULL%1,                  Valerie H.%Mika%NULL%1,                  Tyler%Prentiss%NULL%1,                  John%Zervos%NULL%1,                  Donia A.%Dalal%NULL%1,                  Amber M.%Vasquez%NULL%1,                  Sujan C.%Reddy%NULL%0,                  John%Jernigan%NULL%1,                  Paul E.%Kilgore%NULL%1,                  Marcus J.%Zervos%NULL%1,                  Teena%Chopra%NULL%1,                  Carla P.%Bezold%NULL%1,                  Najibah K.%Rehman%NULL%1]'
$ws.Range("I39").Value = '_PMC_CROSSREF'
$ws.Range("E40").Value = '[Janice K%Louie%Janice.louie@sfdph.org%1,                  Hyman M%Scott%NULL%1,                  Amie%DuBois%NULL%1,                  Natalya%Sturtz%NULL%1,                  Wendy%Lu%NULL%1,                  Juliet%Stoltey%NULL%1,                  Godfred%Masinde%NULL%1,                  Stephanie%Cohen%NULL%1,                  Darpun%Sachdev%NULL%1,                  Susan%Philip%NULL%1,                  Naveena%Bobba%NULL%1,                  Tomas%Aragon%NULL%1,                  NULL%NULL%NULL%0]'
$ws.Range("I40").Value = '_PMC_CROSSREF'
$ws.Range("E41").Value = '[Scott A%Goldberg%sagoldberg@bwh.harvard.edu%1,                  Jochen%Lennerz%NULL%1,                  Michael%Klompas%NULL%1,                  Eden%Mark%NULL%1,                  Virginia M%Pierce%NULL%1,                  Ryan W%Thompson%NULL%1,                  Charles T%Pu%NULL%1,                  Lauren L%Ritterhouse%NULL%1,                  Anand%Dighe%NULL%1,                  Eric S%Rosenberg%NULL%1,                  David C%Grabowski%NULL%1]'
$ws.Range("I41").Value = '_PMC_CROSSREF'
$ws.Range("E42").Value = '[Amy V.%Dora%NULL%1,                  Alexander%Winnett%NULL%1,                  Lauren P.%Jatt%NULL%1,                  Kusha%Davar%NULL%1,                  Mika%Watanabe%NULL%1,                  Linda%Sohn%NULL%1,                  Hannah S.%Kern%NULL%1,                  Christopher J.%Graber%NULL%1,                  Matthew B.%Goetz%NULL%1]'
$ws.Range("I42").Value = '_PMC_CROSSREF'
$ws.Range("E43").Value = '[Benjamin F.%Bigelow%NULL%1,                  Olive%Tang%NULL%1,                  Gregory R.%Toci%NULL%1,                  Norberth%Stracker%NULL%1,                  Fatima%Sheikh%NULL%1,                  Kara M.%Jacobs Slifka%NULL%1,                  Shannon A.%Novosad%NULL%1,                  John A.%Jernigan%NULL%0,                  Sujan C.%Reddy%NULL%0,                  Morgan J.%Katz%NULL%1]'
$ws.Range("I43").Value = '_PMC_CROSSREF'
$ws.Range("E44").Value = '[Matt%Feaster%NULL%1,                  Ying-Ying%Goh%NULL%1]'
$ws.Range("I44").Value = '_PMC_CROSSREF'
$ws.Range("E45").Value = '[William R.%Mills%NULL%1,                  Janet M.%Buccola%NULL%1,                  Susan%Sender%NULL%1,                  Joseph%Lichtefeld%NULL%1,                  Nicholas%Romano%NULL%1,                  Karen%Reynolds%NULL%1,                  Melissa%Price%NULL%1,                  Jennifer%Phipps%NULL%1,                  Leigh%White%NULL%1,                  Shauen%Howard%NULL%1]'
$ws.Range("I45").Value = '_PMC_elsevier_CROSSREF'
$ws.Range("E46").Value = '[Drew A.%Harris%xref no email%1,            Laurie%Archbald-Pannone%xref no email%1,            Jasveen%Kaur%xref no email%1,            David%Cattell-Gordon%xref no email%1,            Karen S.%Rheuban%xref no email%1,            Rachel L.%Ombres%xref no email%1,            Kimberly%Albero%xref no email%1,            Rebecca%Steele%xref no email%1,            Taison D.%Bell%xref no email%1,            Justin B.%Mutter%xref no email%1]'
$ws.Range("I46").Value = '_CROSSREF'
$ws.Range("E47").Value = '[Carl D.%Shrader%xref no email%1,            Shauna%Assadzandi%xref no email%1,            Courtney S.%Pilkerton%xref no email%1,            Amie M.%Ashcraft%xref no email%1]'
$ws.Range("I47").Value = '_CROSSREF'
